{"js": "// 1. \"... a zombie's limb fall off.\" -> \"... a zombie's limb fall off to reduce repetition.\"\n//    (the trailing period is dropped and replaced with a new clause)\nconst weaponLegResults = context.document.body.search(\"a zombie\\u2019s limb fall off.\", { matchCase: true });\nweaponLegResults.load(\"text\");\nawait context.sync();\n\nif (weaponLegResults.items.length > 0) {\n  weaponLegResults.items[0].insertText(\n    \"a zombie\\u2019s limb fall off to reduce repetition.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2. Remove the whole \"BiteAction\" section: the bold heading paragraph, the\n//    paragraph explaining it, and the blank spacer paragraph that follows it.\n//    (The design no longer uses a separate BiteAction class.)\nconst biteHeadingResults = context.document.body.search(\"BiteAction\", { matchCase: true });\nbiteHeadingResults.load(\"text\");\nawait context.sync();\n\nif (biteHeadingResults.items.length > 0) {\n  const headingPara = biteHeadingResults.items[0].paragraphs.getFirst();\n  const bodyPara = headingPara.getNext();\n  const spacerPara = bodyPara.getNext();\n\n  headingPara.delete();\n  bodyPara.delete();\n  spacerPara.delete();\n  await context.sync();\n}\n\n// 3. Update the sentence about the Zombie class' getIntrinsicWeapon method.\n//    3a. Capitalise \"intrinsicweapon\" -> \"Intrinsicweapon\".\nconst intrinsicWordResults = context.document.body.search(\"intrinsicweapon\", { matchCase: true });\nintrinsicWordResults.load(\"text\");\nawait context.sync();\n\nif (intrinsicWordResults.items.length > 0) {\n  intrinsicWordResults.items[0].insertText(\"Intrinsicweapon\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n//    3b. \"... class instance is enough for now.\" ->\n//        \"... class instance is enough for now, to reduce the number of\n//        unnecessary class.\"\nconst instanceResults = context.document.body.search(\n  \"class instance is enough for now.\",\n  { matchCase: true }\n);\ninstanceResults.load(\"text\");\nawait context.sync();\n\nif (instanceResults.items.length > 0) {\n  instanceResults.items[0].insertText(\n    \"class instance is enough for now, to reduce the number of unnecessary class.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"... a zombie's limb fall off.\" -> \"... a zombie's limb fall off to reduce repetition.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Execute(\n    \"limb fall off.\",   # FindText\n    $false,             # MatchCase\n    $false,             # MatchWholeWord\n    $false,             # MatchWildcards\n    $false,             # MatchSoundsLike\n    $false,             # MatchAllWordForms\n    $true,              # Forward\n    1,                  # Wrap (wdFindContinue)\n    $false,             # Format\n    \"limb fall off to reduce repetition.\",  # ReplaceWith\n    2                   # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2. Remove the whole \"BiteAction\" section: the bold heading paragraph, the\n#    paragraph explaining it, and the blank spacer paragraph that follows it.\n#    (The design no longer uses a separate BiteAction class.)\n$paras = $d.Paragraphs\n$biteHeadingIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq \"BiteAction\") {\n        $biteHeadingIndex = $i\n        break\n    }\n}\n\nif ($biteHeadingIndex -gt 0) {\n    $headingPara = $paras.Item($biteHeadingIndex)\n    $spacerPara  = $paras.Item($biteHeadingIndex + 2)\n    $delRange = $d.Range($headingPara.Range.Start, $spacerPara.Range.End)\n    $delRange.Delete()\n}\n\n# 3. Update the sentence about the Zombie class' getIntrinsicWeapon method.\n# 3a. Capitalise \"intrinsicweapon\" -> \"Intrinsicweapon\". MatchCase must be\n#     $true here, otherwise this would also match \"IntrinsicWeapon\" inside\n#     the earlier \"getIntrinsicWeapon\" in the same paragraph.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Execute(\n    \"intrinsicweapon\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Intrinsicweapon\",\n    2\n) | Out-Null\n\n# 3b. \"... class instance is enough for now.\" ->\n#     \"... class instance is enough for now, to reduce the number of\n#     unnecessary class.\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Execute(\n    \"class instance is enough for now.\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"class instance is enough for now, to reduce the number of unnecessary class.\",\n    2\n) | Out-Null\n"}
